# Add a new "2021" column (column R) to the yearly data table on the
# active sheet. Column R did not exist before (data ran through Q/2020),
# so insert a fresh column R; Excel's column-insert duplicates the
# formatting of the column immediately to its left (Q), which is exactly
# the formatting the new 2021 column should carry. Then fill in the 2021
# values on top of that formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(18).Insert()   # column R

# row -> 2021 value (column R). Row 3 is the year header; rows 4-33 are data.
$values = [ordered]@{
    3  = 2021
    4  = 12.6
    5  = 17.9
    6  = 7.3
    7  = 12.6
    8  = 19.5
    9  = 5.5
    10 = 10.3
    11 = 12.3
    12 = 8.2
    13 = 24.8
    14 = 33.1
    15 = 16.6
    16 = 23.9
    17 = 29.5
    18 = 18.1
    19 = 9.6
    20 = 14.8
    21 = 4.3
    22 = 12.1
    23 = 18.2
    24 = 5.9
    25 = 17.3
    26 = 27.6
    27 = 7.4
    28 = 7.8
    29 = 10.4
    30 = 5.6
    31 = 6.7
    32 = 10.7
    33 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 18).Value = $values[$row]
}

# Match the author's recorded selection after the edit.
$ws.Range("S4").Select() | Out-Null
